$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the password value in B2 (the "@AM@TGA1$29TN" corp test credential)
# while keeping the cell's existing style/formatting.
$ws.Range("B2").ClearContents()

# Reset the active cell/selection back to A1 (the stray "B6" selection from
# the prior interactive edit no longer applies once B2's content is gone).
$ws.Range("A1").Select()
